$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A6').Value = 'Remis'
$ws.Range('A8').Value = 'Radomiak Radom'
$ws.Range('A9').Value = 'Warta Poznań'
$ws.Range('A10').Value = 'Lechia Gdańsk'
$ws.Range('A12').Value = 'Raków Częstochowa'
$ws.Range('A14').Value = 'Remis'
$ws.Range('A17').Value = 'Remis'
$ws.Range('A20').Value = 'Legia Warszawa'
$ws.Range('A21').Value = 'Zagłębie Lubin'
$ws.Range('A24').Value = 'Warta Poznań'
$ws.Range('A26').Value = 'Remis'
$ws.Range('A27').Value = 'Raków Częstochowa'
$ws.Range('A28').Value = 'Lechia Gdańsk'
$ws.Range('A29').Value = 'Radomiak Radom'
$ws.Range('A30').Value = 'Lech Poznań'
$ws.Range('A32').Value = 'Remis'
$ws.Range('A37').Value = 'Remis'
$ws.Range('A42').Value = 'Remis'
$ws.Range('A44').Value = 'Lechia Gdańsk'
$ws.Range('A45').Value = 'Remis'
$ws.Range('A49').Value = 'Radomiak Radom'
$ws.Range('A50').Value = 'Remis'
$ws.Range('A52').Value = 'Pogoń Szczecin'
$ws.Range('A53').Value = 'Widzew Łódź'
$ws.Range('A54').Value = 'Remis'
$ws.Range('A55').Value = 'Remis'
$ws.Range('A56').Value = 'Cracovia'
$ws.Range('A58').Value = 'Lechia Gdańsk'
$ws.Range('A61').Value = 'Radomiak Radom'
$ws.Range('A63').Value = 'Remis'
$ws.Range('A65').Value = 'Raków Częstochowa'
$ws.Range('A66').Value = 'Miedź Legnica'
$ws.Range('A67').Value = 'Zagłębie Lubin'
$ws.Range('A72').Value = 'Śląsk Wrocław'
$ws.Range('A76').Value = 'Piast Gliwice'
$ws.Range('A78').Value = 'Remis'
$ws.Range('A79').Value = 'Legia Warszawa'
$ws.Range('A80').Value = 'Remis'
$ws.Range('A81').Value = 'Cracovia'
$ws.Range('A82').Value = 'Lechia Gdańsk'
$ws.Range('A84').Value = 'Piast Gliwice'
$ws.Range('A85').Value = 'Zagłębie Lubin'
$ws.Range('A86').Value = 'Korona Kielce'
$ws.Range('A94').Value = 'Miedź Legnica'
$ws.Range('A96').Value = 'Lechia Gdańsk'
$ws.Range('A97').Value = 'Radomiak Radom'
$ws.Range('A99').Value = 'Wisła Płock'
$ws.Range('A100').Value = 'Warta Poznań'
$ws.Range('A101').Value = 'Remis'
$ws.Range('A102').Value = 'Widzew Łódź'
$ws.Range('A107').Value = 'Miedź Legnica'
$ws.Range('A108').Value = 'Remis'
$ws.Range('A109').Value = 'Remis'
$ws.Range('A110').Value = 'Korona Kielce'
$ws.Range('A112').Value = 'Cracovia'
$ws.Range('A114').Value = 'Pogoń Szczecin'
$ws.Range('A115').Value = 'Radomiak Radom'
$ws.Range('A117').Value = 'Zagłębie Lubin'
$ws.Range('A118').Value = 'Wisła Płock'
$ws.Range('A121').Value = 'Zagłębie Lubin'
$ws.Range('A123').Value = 'Remis'
$ws.Range('A124').Value = 'Stal Mielec'
$ws.Range('A125').Value = 'Górnik Zabrze'
$ws.Range('A126').Value = 'Miedź Legnica'
$ws.Range('A128').Value = 'Jagielonia Białystok'
$ws.Range('A129').Value = 'Cracovia'
$ws.Range('A130').Value = 'Remis'
$ws.Range('A131').Value = 'Lech Poznań'
$ws.Range('A132').Value = 'Remis'
$ws.Range('A137').Value = 'Cracovia'
$ws.Range('A138').Value = 'Piast Gliwice'
$ws.Range('A140').Value = 'Legia Warszawa'
$ws.Range('A141').Value = 'Śląsk Wrocław'
$ws.Range('A142').Value = 'Remis'
$ws.Range('A145').Value = 'Remis'
$ws.Range('A146').Value = 'Remis'
$ws.Range('A147').Value = 'Raków Częstochowa'
$ws.Range('A149').Value = 'Lechia Gdańsk'
$ws.Range('A150').Value = 'Górnik Zabrze'
$ws.Range('A156').Value = 'Jagielonia Białystok'
$ws.Range('A158').Value = 'Korona Kielce'
$ws.Range('A163').Value = 'Remis'
$ws.Range('A164').Value = 'Legia Warszawa'
$ws.Range('A165').Value = 'Cracovia'
$ws.Range('A167').Value = 'Remis'
$ws.Range('A170').Value = 'Raków Częstochowa'
$ws.Range('A171').Value = 'Remis'
$ws.Range('A172').Value = 'Remis'
$ws.Range('A176').Value = 'Legia Warszawa'
$ws.Range('A177').Value = 'Górnik Zabrze'
$ws.Range('A182').Value = 'Remis'
$ws.Range('A183').Value = 'Legia Warszawa'
$ws.Range('A185').Value = 'Remis'
$ws.Range('A188').Value = 'Jagielonia Białystok'
$ws.Range('A193').Value = 'Miedź Legnica'
$ws.Range('A194').Value = 'Lechia Gdańsk'
$ws.Range('A197').Value = 'Warta Poznań'
$ws.Range('A198').Value = 'Wisła Płock'
$ws.Range('A200').Value = 'Śląsk Wrocław'
$ws.Range('A201').Value = 'Remis'
$ws.Range('A203').Value = 'Remis'
$ws.Range('A204').Value = 'Górnik Zabrze'
$ws.Range('A206').Value = 'Remis'
$ws.Range('A207').Value = 'Remis'
$ws.Range('A209').Value = 'Lech Poznań'
$ws.Range('A212').Value = 'Radomiak Radom'
$ws.Range('A216').Value = 'Remis'
$ws.Range('A218').Value = 'Jagielonia Białystok'
$ws.Range('A219').Value = 'Piast Gliwice'
$ws.Range('A220').Value = 'Wisła Płock'
$ws.Range('A227').Value = 'Widzew Łódź'
$ws.Range('A229').Value = 'Warta Poznań'
$ws.Range('A230').Value = 'Remis'
$ws.Range('A231').Value = 'Remis'
$ws.Range('A235').Value = 'Radomiak Radom'
$ws.Range('A236').Value = 'Lechia Gdańsk'
$ws.Range('A242').Value = 'Remis'
$ws.Range('A243').Value = 'Remis'
$ws.Range('A244').Value = 'Śląsk Wrocław'
$ws.Range('A245').Value = 'Remis'
$ws.Range('A249').Value = 'Remis'
$ws.Range('A250').Value = 'Remis'
$ws.Range('A252').Value = 'Remis'
$ws.Range('A253').Value = 'Warta Poznań'
$ws.Range('A254').Value = 'Remis'
$ws.Range('A255').Value = 'Lechia Gdańsk'
$ws.Range('A256').Value = 'Raków Częstochowa'
$ws.Range('A257').Value = 'Remis'
$ws.Range('A258').Value = 'Remis'
$ws.Range('A260').Value = 'Remis'
$ws.Range('A261').Value = 'Widzew Łódź'
$ws.Range('A264').Value = 'Pogoń Szczecin'
$ws.Range('A269').Value = 'Lechia Gdańsk'
$ws.Range('A270').Value = 'Remis'
$ws.Range('A274').Value = 'Remis'
$ws.Range('A276').Value = 'Remis'
$ws.Range('A277').Value = 'Remis'
$ws.Range('A278').Value = 'Pogoń Szczecin'
$ws.Range('A280').Value = 'Remis'
$ws.Range('A287').Value = 'Warta Poznań'
$ws.Range('A288').Value = 'Górnik Zabrze'
$ws.Range('A289').Value = 'Śląsk Wrocław'
$ws.Range('A291').Value = 'Stal Mielec'
$ws.Range('A293').Value = 'Legia Warszawa'
$ws.Range('A294').Value = 'Pogoń Szczecin'
$ws.Range('A298').Value = 'Remis'
$ws.Range('A299').Value = 'Remis'
$ws.Range('A302').Value = 'Legia Warszawa'
$ws.Range('A303').Value = 'Górnik Zabrze'
$ws.Range('A305').Value = 'Remis'
$ws.Range('A307').Value = 'Remis'
